$d = $word.ActiveDocument

$replacements = @(
    @("2023-09-23 Saturday", "2023-09-24 Sunday"),
    @("31×38=1178", "76×15=1140"),
    @("50×84=4200", "31×62=1922"),
    @("30×62=1860", "89×15=1335"),
    @("45×92=4140", "46×92=4232"),
    @("17×56=952", "39×48=1872"),
    @("93×20=1860", "28×48=1344"),
    @("22×20=440", "58×18=1044"),
    @("62×41=2542", "82×71=5822"),
    @("28×67=1876", "35×20=700"),
    @("57×32=1824", "12×13=156"),
    @("62×96=5952", "62×98=6076"),
    @("81×57=4617", "56×72=4032"),
    @("99×64=6336", "75×23=1725"),
    @("38×12=456", "96×55=5280"),
    @("77×27=2079", "44×37=1628"),
    @("67×84=5628", "29×23=667"),
    @("99×92=9108", "27×89=2403"),
    @("75×37=2775", "67×91=6097"),
    @("40×21=840", "68×41=2788"),
    @("82×24=1968", "39×57=2223"),
    @("71×82=5822", "43×88=3784"),
    @("93×23=2139", "79×15=1185"),
    @("91×80=7280", "77×57=4389"),
    @("20×17=340", "97×41=3977"),
    @("14×44=616", "37×13=481")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying $($replacements.Count) replacements"
